# New man pages & other updates. New build.
#
# phylip-programs (sheet1): "treedist.exe" (row 36) now has its R-interface
# counterpart documented: "Rtreedist" in column B, done-date in column C
# (matching the style/format already used by the row above it).
#
# The workbook was also re-saved with the first sheet ("phylip-programs")
# as the active/selected tab (instead of "addt'l-functions"), scrolled so
# the newly edited row is in view, with B37 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- add the new "Rtreedist" row-36 entries on the phylip-programs sheet ---
$ws1.Range("B36").Value = "Rtreedist"
$ws1.Range("C36").Value = 41613

# Match C36's date formatting to the identically-shaped cell above it (C35)
# without minting a brand new style: copy formats only.
$ws1.Range("C35").Copy()
$ws1.Range("C36").PasteSpecial(-4122)

# --- make "phylip-programs" the active sheet/tab, with B37 selected ---
$ws1.Activate()
$ws1.Range("B37").Select()
